$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Exclude" -> "Include"
$ws.Range("F1").Value2 = "Include"

# Flip F2:F21 values (0 -> 1, 1 -> 0)
for ($r = 2; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    if ($cell.Value2 -eq 0) {
        $cell.Value2 = 1
    } else {
        $cell.Value2 = 0
    }
}

# Update selection to B11
$ws.Range("B11").Select()
